$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-17 with the new/rotated store data ---

# Row 2: Assaí Guaianases
$ws.Range("A2").Value = "Assaí Guaianases"
$ws.Range("B2").Value = "Estrada Dom João Nery (com Rua Engenheiro Bardot), 4.031"
$ws.Range("C2").Value = "São Paulo"
$ws.Range("D2").Value = "SP"
$ws.Range("E2").Value = -23.526843
$ws.Range("F2").Value = -46.396465999999997
$ws.Range("G2").Value = "Nova"
$ws.Range("H2").Value = "3T2022"

# Row 3: Assaí Guaianases Estação
$ws.Range("A3").Value = "Assaí Guaianases Estação"
$ws.Range("B3").Value = "Estrada Itaquera Guaianases"
$ws.Range("C3").Value = "São Paulo"
$ws.Range("D3").Value = "SP"
$ws.Range("E3").Value = -23.542258
$ws.Range("F3").Value = -46.423996000000002
$ws.Range("G3").Value = "Antigas"
$ws.Range("H3").ClearContents()

# Row 4: Assaí Interlagos
$ws.Range("A4").Value = "Assaí Interlagos"
$ws.Range("B4").Value = "Av Sarg Geraldo Santa'ana"
$ws.Range("C4").Value = "São Paulo"
$ws.Range("D4").Value = "SP"
$ws.Range("E4").Value = -23.663383
$ws.Range("F4").Value = -46.680224000000003
$ws.Range("G4").Value = "Nova"
$ws.Range("H4").Value = "3T2022"

# Row 5: Assaí Itapevi
$ws.Range("A5").Value = "Assaí Itapevi"
$ws.Range("B5").Value = "Rod. Engenheiro Renê Benedito da Silva, 977"
$ws.Range("C5").Value = "São Paulo"
$ws.Range("D5").Value = "SP"
$ws.Range("E5").Value = -23.546430999999998
$ws.Range("F5").Value = -46.946899999999999
$ws.Range("G5").Value = "Antigas"
$ws.Range("H5").ClearContents()

# Row 6: Assaí Itaquera
$ws.Range("A6").Value = "Assaí Itaquera"
$ws.Range("B6").Value = "Avenida Sylvio Torres, 190"
$ws.Range("C6").Value = "São Paulo"
$ws.Range("D6").Value = "SP"
$ws.Range("E6").Value = -23.553196
$ws.Range("F6").Value = -46.488427999999999
$ws.Range("G6").Value = "Antigas"
$ws.Range("H6").ClearContents()

# Row 7: Assaí Jaçanã
$ws.Range("A7").Value = "Assaí Jaçanã"
$ws.Range("B7").Value = "Avenida Luís Stamatis, 35"
$ws.Range("C7").Value = "São Paulo"
$ws.Range("D7").Value = "SP"
$ws.Range("E7").Value = -23.467191
$ws.Range("F7").Value = -46.584729000000003
$ws.Range("G7").Value = "Antigas"
$ws.Range("H7").ClearContents()

# Row 8: Assaí Jabaquara
$ws.Range("A8").Value = "Assaí Jabaquara"
$ws.Range("B8").Value = "Rua Taquaruçu, 79"
$ws.Range("C8").Value = "São Paulo"
$ws.Range("D8").Value = "SP"
$ws.Range("E8").Value = -23.641707
$ws.Range("F8").Value = -46.644190999999999
$ws.Range("G8").Value = "Antigas"
$ws.Range("H8").ClearContents()

# Row 9: Assaí Jacu Pêssego
$ws.Range("A9").Value = "Assaí Jacu Pêssego"
$ws.Range("B9").Value = "Avenida Jacu Pêssego, 750"
$ws.Range("C9").Value = "São Paulo"
$ws.Range("D9").Value = "SP"
$ws.Range("E9").Value = -23.565317
$ws.Range("F9").Value = -46.446522000000002
$ws.Range("G9").Value = "Antigas"
$ws.Range("H9").ClearContents()

# Row 10: Assaí Jaguaré
$ws.Range("A10").Value = "Assaí Jaguaré"
$ws.Range("B10").Value = "Avenida Jaguaré, 925"
$ws.Range("C10").Value = "São Paulo"
$ws.Range("D10").Value = "SP"
$ws.Range("E10").Value = -23.552852999999999
$ws.Range("F10").Value = -46.742243999999999
$ws.Range("G10").Value = "Antigas"
$ws.Range("H10").ClearContents()

# Row 11: Assaí Jaguaré Corifeu
$ws.Range("A11").Value = "Assaí Jaguaré Corifeu"
$ws.Range("B11").Value = "Av Corifeu de Azevedo Marques, Jaguaré"
$ws.Range("C11").Value = "São Paulo"
$ws.Range("D11").Value = "SP"
$ws.Range("E11").Value = -23.556875000000002
$ws.Range("F11").Value = -46.748359000000001
$ws.Range("G11").Value = "Nova"
$ws.Range("H11").Value = "4T2022"

# Row 12: Assaí Jaraguá/Taipas
$ws.Range("A12").Value = "Assaí Jaraguá/Taipas"
$ws.Range("B12").Value = "Avenida Raimundo Pereira de Magalhães, 10.535"
$ws.Range("C12").Value = "São Paulo"
$ws.Range("D12").Value = "SP"
$ws.Range("E12").Value = -23.449862
$ws.Range("F12").Value = -46.722527999999997
$ws.Range("G12").Value = "Antigas"
$ws.Range("H12").ClearContents()

# Row 13: Assaí Nações Unidas (new)
$ws.Range("A13").Value = "Assaí Nações Unidas"
$ws.Range("B13").Value = "Av. das Nações Unidas"
$ws.Range("C13").Value = "São Paulo"
$ws.Range("D13").Value = "SP"
$ws.Range("E13").Value = -23.678526999999999
$ws.Range("F13").Value = -46.695574999999998
$ws.Range("G13").Value = "Antigas"
$ws.Range("H13").ClearContents()

# Row 14: Assaí Nordestina (new)
$ws.Range("A14").Value = "Assaí Nordestina"
$ws.Range("B14").Value = "Avenida Nordestina, 3.077"
$ws.Range("C14").Value = "São Paulo"
$ws.Range("D14").Value = "SP"
$ws.Range("E14").Value = -23.510605000000002
$ws.Range("F14").Value = -46.430522000000003
$ws.Range("G14").Value = "Antigas"
$ws.Range("H14").ClearContents()

# Row 15: Assaí Penha - Marginal Tietê (new)
$ws.Range("A15").Value = "Assaí Penha - Marginal Tietê"
$ws.Range("B15").Value = "Av. Condessa Elizabeth de Robiano, 5500"
$ws.Range("C15").Value = "São Paulo"
$ws.Range("D15").Value = "SP"
$ws.Range("E15").Value = -23.512795000000001
$ws.Range("F15").Value = -46.553893000000002
$ws.Range("G15").Value = "Antigas"
$ws.Range("H15").ClearContents()

# Row 16: Assaí Penha Tiquatira (new)
$ws.Range("A16").Value = "Assaí Penha Tiquatira"
$ws.Range("B16").Value = "Av. São Miguel"
$ws.Range("C16").Value = "São Paulo"
$ws.Range("D16").Value = "SP"
$ws.Range("E16").Value = -23.515934000000001
$ws.Range("F16").Value = -46.518680000000003
$ws.Range("G16").Value = "Nova"
$ws.Range("H16").Value = "4T2022"

# Row 17: Assaí Raposo Tavares (new)
$ws.Range("A17").Value = "Assaí Raposo Tavares"
$ws.Range("B17").Value = "Av. Mal. Fiuza de Castro"
$ws.Range("C17").Value = "São Paulo"
$ws.Range("D17").Value = "SP"
$ws.Range("E17").Value = -23.584363
$ws.Range("F17").Value = -46.747152
$ws.Range("G17").Value = "Nova"
$ws.Range("H17").Value = "4T2022"

# --- Remove the trailing blank rows 36-43 so the table shrinks to A1:H35 ---
$ws.Range("A36:H43").EntireRow.Delete()

# --- Update the active selection to H1 ---
$ws.Range("H1").Select()
